$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.3375586344204518, 1.4941583295314556, -0.8257965719468602, 2.302340791781454)"
$ws.Range("C2").Value = "MIE(7.536251712360564, 5.5803136488842915, -14.837131736427386, 22.205954556824096)"
$ws.Range("D2").Value = "JSU(-1.0195015450698808, 1.2851424364099886, 0.2724282462206998, 3.2260262461841496)"
$ws.Range("E2").Value = "JSU(-1.1522416378847942, 1.0014978470275973, 3.2787213957246646, 3.573841903868784)"
